$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

Replace-Text "Introspektivni pogled na odnos između Hawkinga i prostor/vremenskog kontingenta. Ovaj film izražava " "Introspektivni pogled na odnos između Hawkinga i prostor/vremenskog kontingenta. ## Ovaj film izražava "
Replace-Text "odnos prema Einsteinovoj teoriji opće relativnosti. Film je " "odnos prema Einsteinovoj teoriji opće relativnosti. ## Film je "
Replace-Text "čovjeka (Hawking) kao i njegov rad (Crne rupe). Intervjui sa svojom obitelji su malo " "čovjeka (Hawking) kao i njegov rad (Crne rupe). ## Intervjui sa svojom obitelji su malo "
Replace-Text " njegove teorije i ideje. Philip Glass " " njegove teorije i ideje. ## Philip Glass "
Replace-Text " filmu. Samo je" " filmu. ## Samo je"
Replace-Text "melodije (Jean Michel Jarre). Sveukupno bih " "melodije (Jean Michel Jarre). ## Sveukupno bih "
Replace-Text " dugo... dugo vremena... " " dugo... dugo vremena... ## "
Replace-Text "sam ga sinoć i htio otići nakon 20 minuta... Keira Knightley" "sam ga sinoć i htio otići nakon 20 minuta... ## Keira Knightley"
Replace-Text "karizmu ispuniti ulogu... " "karizmu ispuniti ulogu...  ##"
Replace-Text " Je li ikad imala satove glume? Sudeći po " " Je li ikad imala satove glume? ## Sudeći po "
Replace-Text "u bliskoj budućnosti... " "u bliskoj budućnosti... ## "
Replace-Text ".. ako " ".. ## ako "
Replace-Text "Zahvaljujući drugim recenzentima koji su me usmjerili na ovaj proizvod kad mi je rečeno da sam anemična. " "Zahvaljujući drugim recenzentima koji su me usmjerili na ovaj proizvod kad mi je rečeno da sam anemična. ## "
Replace-Text "anemija je nestala. Dobar proizvod. Jednostavno " "anemija je nestala. ## Dobar proizvod. ## Jednostavno "
Replace-Text "Ovo je jedan od mojih omiljenih deserta i brzo se topi u ustima. Ova marka je dobra i isporučena dobro zapakirana. Svatko bi trebao probati " "Ovo je jedan od mojih omiljenih deserta i brzo se topi u ustima. ## Ova marka je dobra i isporučena dobro zapakirana. ## Svatko bi trebao probati "
Replace-Text "jednom. Cijena " "jednom. ## Cijena "
Replace-Text "Ovo je fantastična zagonetka/poklon za mlade i stare." "Ovo je fantastična zagonetka/poklon za mlade i stare. ##"
Replace-Text "zajedno na širok broj načina. " "zajedno na širok broj načina. ## "
Replace-Text ". U usporedbi s većinom drugih, jedina razlika ovdje je " ". ## U usporedbi s većinom drugih, jedina razlika ovdje je "
Replace-Text " žensko. " " žensko. ## "
Replace-Text "je ista. Akcijske scene nisu zanimljive. Specijalni efekti su " "je ista. ## Akcijske scene nisu zanimljive. ## Specijalni efekti su "
Replace-Text ". Bez obzira na " ". ## Bez obzira na "
Replace-Text " naginjati na jednu stranu. " " naginjati na jednu stranu. ## "
Replace-Text "ili raditi trikove. Imam " "ili raditi trikove. ## Imam "
Replace-Text "prekratka. Nemojte ih dobiti ako udarate tešku " "prekratka. ## Nemojte ih dobiti ako udarate tešku "
Replace-Text ". Oni jednostavno neće štititi/podupirati vaše zglobove ili " ". ## Oni jednostavno neće štititi/podupirati vaše zglobove ili "

# Insert a new empty paragraph after the paragraph that ends with
# "...vaše zglobove ili zglobove." and before the final (already existing)
# empty paragraph.
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd()
    if ($t.EndsWith("štititi/podupirati vaše zglobove ili zglobove.")) {
        $targetPara = $p
    }
}
if ($targetPara -ne $null) {
    $targetPara.Range.InsertParagraphAfter()
} else {
    Write-Host "Target paragraph for insertion not found"
}
